$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Monday")

# New Monday speeches from Cision (rows 10-18)
$rows = @(
    @("Luttrell",  "Marcus",  "rnc.luttrell.txt"),
    @("Smith",     "Pat",     "rnc.smithpat.txt"),
    @("Geist",     "Mark",    "rnc.geist.txt"),
    @("Sabato",    "Antonio", "rnc.sabato.txt"),
    @("McCaul",    "Mike",    "rnc.mccaul.txt"),
    @("Clarke",    "David",   "rnc.clarke.txt"),
    @("Cotton",    "Tim",     "rnc.cotton.txt"),
    @("Beardsley", "Jason",   "rnc.beardsley.txt"),
    @("Zinke",     "Ryan",    "rnc.zinke.txt")
)

$r = 10
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "Monday"
    $ws.Cells.Item($r, 5).Value = "speech"
    $ws.Cells.Item($r, 6).Value = "Cision"
    $r = $r + 1
}

$ws.Range("F19").Select()

# Tuesday becomes the active/selected tab
$tue = $wb.Worksheets.Item("Tuesday")
$tue.Activate()
